# Fruta / hortaliza, semanal
# Insert a new weekly record at row 15 (pushing the existing rows 15-26
# down to 16-27) and populate the new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 15:26 down to 16:27, carrying their formatting along.
$ws.Rows("15:15").Insert()

# New row 15 - constant/shared columns match the rest of the series.
$ws.Range("A15").Value2 = 10
$ws.Range("B15").Value2 = "Vega Modelo de Temuco"
$ws.Range("C15").Value2 = "La Araucanía"
$ws.Range("D15").Value2 = 45134
$ws.Range("E15").Value2 = 9
$ws.Range("F15").Value2 = 100112036
$ws.Range("G15").Value2 = "Caigua"
$ws.Range("H15").Value2 = "Sin especificar"
$ws.Range("I15").Value2 = "Primera"
$ws.Range("J15").Value2 = 5
$ws.Range("K15").Value2 = 20000
$ws.Range("L15").Value2 = 20000
$ws.Range("M15").Value2 = 20000
$ws.Range("N15").Value2 = "$/caja 15 kilos"
$ws.Range("O15").Value2 = "Región de Arica y Parinacota"
$ws.Range("P15").Value2 = 1333
$ws.Range("Q15").Value2 = 15
$ws.Range("R15").Value2 = "Hortaliza"
